$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headers = @(
    "origen",
    "periodo",
    "A",
    "A1",
    "A2",
    "A3",
    "A4",
    "A5",
    "A6",
    "A7",
    "A8",
    "B",
    "B1",
    "B2",
    "B3",
    "B4",
    "B5",
    "B6",
    "trabajadores_unicos",
    "empleadores_unicos",
    "planillas_unicas"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(2, $i + 1).Value = $headers[$i]
}

$ws.Range("A2:U2").Select()
